$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 981.2143
$ws.Cells.Item(28, 9).Value = 751.4167
$ws.Cells.Item(28, 10).Value = 2360
$ws.Cells.Item(28, 11).Value = 751.4167
$ws.Cells.Item(28, 12).Value = 2360
$ws.Cells.Item(28, 13).Value = -266.4167
$ws.Cells.Item(28, 14).Value = -3330

$ws.Cells.Item(107, 8).Value = 3767.682
$ws.Cells.Item(107, 9).Value = 3478.4736
$ws.Cells.Item(107, 11).Value = 3478.4736
$ws.Cells.Item(107, 13).Value = -1558.4736

$ws.Cells.Item(112, 8).Value = 1780.625
$ws.Cells.Item(112, 10).Value = 1678.2142
$ws.Cells.Item(112, 12).Value = 5034.642599999999
$ws.Cells.Item(112, 14).Value = -7250.642599999999

$ws.Cells.Item(137, 8).Value = 6259703.5
$ws.Cells.Item(137, 9).Value = 20003480
$ws.Cells.Item(137, 10).Value = 12533.137
$ws.Cells.Item(137, 11).Value = 60010440
$ws.Cells.Item(137, 12).Value = 37599.411
$ws.Cells.Item(137, 13).Value = -60007890
$ws.Cells.Item(137, 14).Value = -42699.411

$ws.Cells.Item(138, 8).Value = 6604.4893
$ws.Cells.Item(138, 9).Value = 6184.643
$ws.Cells.Item(138, 10).Value = 6782.606
$ws.Cells.Item(138, 11).Value = 18553.929
$ws.Cells.Item(138, 12).Value = 20347.818
$ws.Cells.Item(138, 13).Value = -13413.929
$ws.Cells.Item(138, 14).Value = -30627.818

$ws.Cells.Item(141, 8).Value = 5746.7896
$ws.Cells.Item(141, 9).Value = 2608.7273
$ws.Cells.Item(141, 11).Value = 7826.1819
$ws.Cells.Item(141, 13).Value = -2646.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 995.4
$ws.Cells.Item(2, 9).Value = 995.4
$ws.Cells.Item(2, 11).Value = 995.4
$ws.Cells.Item(2, 13).Value = -882.4

$ws.Cells.Item(32, 8).Value = 706710.75
$ws.Cells.Item(32, 9).Value = 725666.1
$ws.Cells.Item(32, 10).Value = 52750
$ws.Cells.Item(32, 11).Value = 725666.1
$ws.Cells.Item(32, 12).Value = 52750
$ws.Cells.Item(32, 13).Value = -725379.1
$ws.Cells.Item(32, 14).Value = -53324

$ws.Cells.Item(61, 8).Value = 4303419
$ws.Cells.Item(61, 9).Value = 3335397.5
$ws.Cells.Item(61, 10).Value = 8336841
$ws.Cells.Item(61, 11).Value = 3335397.5
$ws.Cells.Item(61, 12).Value = 8336841
$ws.Cells.Item(61, 13).Value = -3335185.5
$ws.Cells.Item(61, 14).Value = -8337265

$ws.Cells.Item(74, 8).Value = 2082168.8
$ws.Cells.Item(74, 9).Value = 2422762
$ws.Cells.Item(74, 11).Value = 2422762
$ws.Cells.Item(74, 13).Value = -2421888

$ws.Cells.Item(77, 8).Value = 2082168.8
$ws.Cells.Item(77, 9).Value = 2422762
$ws.Cells.Item(77, 11).Value = 12113810
$ws.Cells.Item(77, 13).Value = -12109442

$ws.Cells.Item(116, 8).Value = 995.4
$ws.Cells.Item(116, 9).Value = 995.4
$ws.Cells.Item(116, 11).Value = 995.4
$ws.Cells.Item(116, 13).Value = 1298.6

$ws.Cells.Item(132, 8).Value = 4008.1865
$ws.Cells.Item(132, 9).Value = 2818.1177
$ws.Cells.Item(132, 10).Value = 5626.68
$ws.Cells.Item(132, 11).Value = 8454.3531
$ws.Cells.Item(132, 12).Value = 16880.04
$ws.Cells.Item(132, 13).Value = -5924.3531
$ws.Cells.Item(132, 14).Value = -21940.04

$ws.Cells.Item(136, 8).Value = 4303419
$ws.Cells.Item(136, 9).Value = 3335397.5
$ws.Cells.Item(136, 10).Value = 8336841
$ws.Cells.Item(136, 11).Value = 10006192.5
$ws.Cells.Item(136, 12).Value = 25010523
$ws.Cells.Item(136, 13).Value = -10003642.5
$ws.Cells.Item(136, 14).Value = -25015623

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 995.4
$ws.Cells.Item(3, 9).Value = 995.4
$ws.Cells.Item(3, 11).Value = 995.4
$ws.Cells.Item(3, 13).Value = -881.4

$ws.Cells.Item(36, 8).Value = 1342.375
$ws.Cells.Item(36, 9).Value = 1342.375
$ws.Cells.Item(36, 11).Value = 1342.375
$ws.Cells.Item(36, 13).Value = -808.375

$ws.Cells.Item(86, 8).Value = 1572.56
$ws.Cells.Item(86, 9).Value = 1395.8125
$ws.Cells.Item(86, 10).Value = 1886.7778
$ws.Cells.Item(86, 11).Value = 1395.8125
$ws.Cells.Item(86, 12).Value = 1886.7778
$ws.Cells.Item(86, 13).Value = -272.8125
$ws.Cells.Item(86, 14).Value = -4132.7778

$ws.Cells.Item(89, 8).Value = 1572.56
$ws.Cells.Item(89, 9).Value = 1395.8125
$ws.Cells.Item(89, 10).Value = 1886.7778
$ws.Cells.Item(89, 11).Value = 6979.0625
$ws.Cells.Item(89, 12).Value = 9433.889000000001
$ws.Cells.Item(89, 13).Value = -1363.0625
$ws.Cells.Item(89, 14).Value = -20665.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 327907.9
$ws.Cells.Item(31, 9).Value = 600754.25
$ws.Cells.Item(31, 11).Value = 600754.25
$ws.Cells.Item(31, 13).Value = -600459.25

$ws.Cells.Item(34, 8).Value = 327907.9
$ws.Cells.Item(34, 9).Value = 600754.25
$ws.Cells.Item(34, 11).Value = 600754.25
$ws.Cells.Item(34, 13).Value = -600552.25

$ws.Cells.Item(58, 8).Value = 5682265
$ws.Cells.Item(58, 9).Value = 15154322
$ws.Cells.Item(58, 10).Value = 1674856.2
$ws.Cells.Item(58, 11).Value = 15154322
$ws.Cells.Item(58, 12).Value = 1674856.2
$ws.Cells.Item(58, 13).Value = -15154119
$ws.Cells.Item(58, 14).Value = -1675262.2

$ws.Cells.Item(86, 8).Value = 43376.22
$ws.Cells.Item(86, 9).Value = 202446.25
$ws.Cells.Item(86, 10).Value = 9887.789000000001
$ws.Cells.Item(86, 11).Value = 202446.25
$ws.Cells.Item(86, 12).Value = 9887.789000000001
$ws.Cells.Item(86, 13).Value = -201323.25
$ws.Cells.Item(86, 14).Value = -12133.789

$ws.Cells.Item(89, 8).Value = 43376.22
$ws.Cells.Item(89, 9).Value = 202446.25
$ws.Cells.Item(89, 10).Value = 9887.789000000001
$ws.Cells.Item(89, 11).Value = 1012231.25
$ws.Cells.Item(89, 12).Value = 49438.94500000001
$ws.Cells.Item(89, 13).Value = -1006615.25
$ws.Cells.Item(89, 14).Value = -60670.94500000001

$ws.Cells.Item(132, 8).Value = 2749.56
$ws.Cells.Item(132, 10).Value = 3465.6667
$ws.Cells.Item(132, 12).Value = 10397.0001
$ws.Cells.Item(132, 14).Value = -15457.0001

$ws.Cells.Item(136, 8).Value = 5682265
$ws.Cells.Item(136, 9).Value = 15154322
$ws.Cells.Item(136, 10).Value = 1674856.2
$ws.Cells.Item(136, 11).Value = 45462966
$ws.Cells.Item(136, 12).Value = 5024568.6
$ws.Cells.Item(136, 13).Value = -45460416
$ws.Cells.Item(136, 14).Value = -5029668.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(61, 8).Value = 73.333336
$ws.Cells.Item(61, 10).Value = 42.5
$ws.Cells.Item(61, 12).Value = 127.5
$ws.Cells.Item(61, 14).Value = -557.5

$ws.Cells.Item(92, 8).Value = 1144.1923
$ws.Cells.Item(92, 10).Value = 1805.4445
$ws.Cells.Item(92, 12).Value = 5416.333500000001
$ws.Cells.Item(92, 14).Value = -7912.333500000001

$ws.Cells.Item(126, 8).Value = 11666.5
$ws.Cells.Item(126, 9).Value = 4999.5
$ws.Cells.Item(126, 11).Value = 14998.5
$ws.Cells.Item(126, 13).Value = -10058.5

$ws.Cells.Item(131, 8).Value = 2510.111
$ws.Cells.Item(131, 10).Value = 5713.7144
$ws.Cells.Item(131, 12).Value = 17141.1432
$ws.Cells.Item(131, 14).Value = -27221.1432

$ws.Cells.Item(133, 8).Value = 3515.05
$ws.Cells.Item(133, 9).Value = 2910.5789
$ws.Cells.Item(133, 11).Value = 8731.736699999999
$ws.Cells.Item(133, 13).Value = -3671.736699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1296.2858
$ws.Cells.Item(80, 9).Value = 1186.625
$ws.Cells.Item(80, 10).Value = 1442.5
$ws.Cells.Item(80, 11).Value = 1186.625
$ws.Cells.Item(80, 12).Value = 1442.5
$ws.Cells.Item(80, 13).Value = -188.625
$ws.Cells.Item(80, 14).Value = -3438.5

$ws.Cells.Item(83, 8).Value = 1296.2858
$ws.Cells.Item(83, 9).Value = 1186.625
$ws.Cells.Item(83, 10).Value = 1442.5
$ws.Cells.Item(83, 11).Value = 5933.125
$ws.Cells.Item(83, 12).Value = 7212.5
$ws.Cells.Item(83, 13).Value = -941.125
$ws.Cells.Item(83, 14).Value = -17196.5

$ws.Cells.Item(102, 8).Value = 2576.818
$ws.Cells.Item(102, 9).Value = 1345.8462
$ws.Cells.Item(102, 11).Value = 1345.8462
$ws.Cells.Item(102, 13).Value = 276.1538

$ws.Cells.Item(107, 8).Value = 900.8
$ws.Cells.Item(107, 9).Value = 897.625
$ws.Cells.Item(107, 10).Value = 906.44446
$ws.Cells.Item(107, 11).Value = 897.625
$ws.Cells.Item(107, 12).Value = 906.44446
$ws.Cells.Item(107, 13).Value = 1022.375
$ws.Cells.Item(107, 14).Value = -4746.44446

$ws.Cells.Item(113, 8).Value = 1805.4615
$ws.Cells.Item(113, 9).Value = 2000.8182
$ws.Cells.Item(113, 10).Value = 731
$ws.Cells.Item(113, 11).Value = 2000.8182
$ws.Cells.Item(113, 12).Value = 731
$ws.Cells.Item(113, 13).Value = 169.1818000000001
$ws.Cells.Item(113, 14).Value = -5071

$ws.Cells.Item(122, 8).Value = 5075377.5
$ws.Cells.Item(122, 9).Value = 10147899
$ws.Cells.Item(122, 11).Value = 30443697
$ws.Cells.Item(122, 13).Value = -30441247

$ws.Cells.Item(123, 8).Value = 73326
$ws.Cells.Item(123, 10).Value = 73326
$ws.Cells.Item(123, 12).Value = 73326
$ws.Cells.Item(123, 14).Value = -78226

$ws.Cells.Item(132, 8).Value = 16108.853
$ws.Cells.Item(132, 9).Value = 13580.25
$ws.Cells.Item(132, 10).Value = 22177.5
$ws.Cells.Item(132, 11).Value = 40740.75
$ws.Cells.Item(132, 12).Value = 66532.5
$ws.Cells.Item(132, 13).Value = -38210.75
$ws.Cells.Item(132, 14).Value = -71592.5

$ws.Cells.Item(136, 8).Value = 69000
$ws.Cells.Item(136, 10).Value = 69000
$ws.Cells.Item(136, 12).Value = 207000
$ws.Cells.Item(136, 14).Value = -212100

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1551.2858
$ws.Cells.Item(16, 9).Value = 1535.3334
$ws.Cells.Item(16, 11).Value = 1535.3334
$ws.Cells.Item(16, 13).Value = -1365.3334

$ws.Cells.Item(40, 8).Value = 5432.154
$ws.Cells.Item(40, 9).Value = 5419.8184
$ws.Cells.Item(40, 10).Value = 5500
$ws.Cells.Item(40, 11).Value = 5419.8184
$ws.Cells.Item(40, 12).Value = 5500
$ws.Cells.Item(40, 13).Value = -5283.8184
$ws.Cells.Item(40, 14).Value = -5772

$ws.Cells.Item(132, 8).Value = 2607552.5
$ws.Cells.Item(132, 9).Value = 4169730
$ws.Cells.Item(132, 11).Value = 12509190
$ws.Cells.Item(132, 13).Value = -12506660

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 31243.4
$ws.Cells.Item(45, 10).Value = 32804.25
$ws.Cells.Item(45, 12).Value = 32804.25
$ws.Cells.Item(45, 14).Value = -33786.25

$ws.Cells.Item(62, 8).Value = 21180.4
$ws.Cells.Item(62, 9).Value = 1966.3334
$ws.Cells.Item(62, 10).Value = 50001.5
$ws.Cells.Item(62, 11).Value = 1966.3334
$ws.Cells.Item(62, 12).Value = 50001.5
$ws.Cells.Item(62, 13).Value = -1342.3334
$ws.Cells.Item(62, 14).Value = -51249.5

$ws.Cells.Item(65, 8).Value = 21180.4
$ws.Cells.Item(65, 9).Value = 1966.3334
$ws.Cells.Item(65, 10).Value = 50001.5
$ws.Cells.Item(65, 11).Value = 9831.666999999999
$ws.Cells.Item(65, 12).Value = 250007.5
$ws.Cells.Item(65, 13).Value = -6711.666999999999
$ws.Cells.Item(65, 14).Value = -256247.5

$ws.Cells.Item(126, 8).Value = 3395.6155
$ws.Cells.Item(126, 9).Value = 3688.3333
$ws.Cells.Item(126, 11).Value = 11064.9999
$ws.Cells.Item(126, 13).Value = -8594.999899999999

$ws.Cells.Item(136, 8).Value = 1768222.4
$ws.Cells.Item(136, 9).Value = 1036379.4
$ws.Cells.Item(136, 11).Value = 3109138.2
$ws.Cells.Item(136, 13).Value = -3106588.2
